$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: version 0.6.4 (values entered B, A, C, D, E, F to match original authoring order)
$ws.Range("B23").Value = "AUTOMATA CELULAR - copia (33)"
$ws.Range("A23").Value = "0.6.4"
$ws.Range("C23").Value = "-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Work on Group selection."
$ws.Range("D23").Value = "-Fixed Greed calc`n-Fixed historic data writing"
$ws.Range("E23").Value = "Python 3.6.1"
$ws.Range("F23").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"

# Row 24: version 0.7.0 (values entered A, B, D, C, E, F to match original authoring order)
$ws.Range("A24").Value = "0.7.0"
$ws.Range("B24").Value = "AUTOMATA CELULAR - copia (34)"
$ws.Range("D24").Value = "-Rework GS.`n-GUI: added Lambda"
$ws.Range("C24").Value = "-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Rework flexibility"
$ws.Range("E24").Value = "Python 3.6.1"
$ws.Range("F24").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"

# Row heights
$ws.Rows.Item(23).RowHeight = 72
$ws.Rows.Item(24).RowHeight = 72

# View: pane frozen scroll position + selection
$ws.Range("C25").Select()
$excel.ActiveWindow.ScrollRow = 22
